# Add v16 Plan runs (2035, 2050) to ModelRuns sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModelRuns")

# ---------------------------------------------------------------------
# 1. Clear the "status" (and, for row 205, "run_offmodel") markers that
#    are no longer current now that newer runs exist.
# ---------------------------------------------------------------------
$ws.Range("I26").Clear()
$ws.Range("I86").Clear()
$ws.Range("I87").Clear()
$ws.Range("I179").Clear()
$ws.Range("I205").Clear()
$ws.Range("J205").Clear()

# ---------------------------------------------------------------------
# 2. Insert a new row for the 2035 Plan v16 run right after the existing
#    2035_TM161_FBP_Plan_15 row (row 206), copying formatting from it.
# ---------------------------------------------------------------------
$ws.Rows.Item(207).Insert()

$ws.Range("A207").Value() = 2035
$ws.Range("B207").Value() = "2035_TM161_FBP_Plan_16"
$ws.Range("C207").Value() = "RTP2025"
$ws.Range("D207").Value() = "FBP"
$ws.Range("F207").Value() = "Update network (v35)"
$ws.Range("G207").Value() = "BAUS_FBP_v08\2035"
$ws.Range("H207").Value() = "PBA50Plus_Final_Blueprint_v65"
$ws.Range("I207").Value() = "current"
$ws.Range("J207").Value() = "FBP"
$ws.Range("K207").Value() = "BlueprintNetworks_v35\net_2035_Blueprint"
$ws.Range("L207").Value() = "model3-b"
$ws.Range("M207").Value() = "https://app.asana.com/1/11860278793487/project/1204085012544660/task/1210118366389838?focus=true"
$ws.Range("N207").Value() = 15.66
$ws.Range("O207").Value() = "na"
$ws.Range("P207").Value() = "na"
$ws.Range("T207").Value() = -0.455
$ws.Range("U207").Value() = 5
$ws.Range("V207").Value() = 27
$ws.Range("W207").Value() = 0
$ws.Range("X207").Value() = 60
$ws.Range("Y207").Value() = "2035 Plan"

# ---------------------------------------------------------------------
# 3. Append a new row for the 2050 Plan v16 run at the very end of the
#    table (after what is now row 246, formerly 2050_TM161_FBP_Plan_15),
#    copying formatting from it.
# ---------------------------------------------------------------------
$ws.Rows.Item(247).Insert()

$ws.Range("A247").Value() = 2050
$ws.Range("B247").Value() = "2050_TM161_FBP_Plan_16"
$ws.Range("C247").Value() = "RTP2025"
$ws.Range("D247").Value() = "FBP"
$ws.Range("F247").Value() = "Update network (v35)"
$ws.Range("G247").Value() = "BAUS_FBP_v08\2050"
$ws.Range("H247").Value() = "PBA50Plus_Final_Blueprint_v65"
$ws.Range("I247").Value() = "current"
$ws.Range("J247").Value() = "FBP"
$ws.Range("K247").Value() = "BlueprintNetworks_v35\net_2050_Blueprint"
$ws.Range("L247").Value() = "model3-c"
$ws.Range("M247").Value() = "https://app.asana.com/1/11860278793487/project/1204085012544660/task/1210118366389840?focus=true"
$ws.Range("N247").Value() = 16.47
$ws.Range("O247").Value() = "na"
$ws.Range("P247").Value() = "na"
$ws.Range("T247").Value() = -0.455
$ws.Range("U247").Value() = 5
$ws.Range("V247").Value() = 34
$ws.Range("W247").Value() = 0
$ws.Range("X247").Value() = 85
$ws.Range("Y247").Value() = "2050 Plan"

# The "2050 Plan" alias marker moves from the previous latest FBP Plan
# run (now row 246) to the new v16 run (row 247).
$ws.Range("Y246").ClearContents()

# ---------------------------------------------------------------------
# 4. Update the frozen-pane view so the newly added rows are visible.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 209
$ws.Range("A247").Select()
